# Annual Vehicle Maint Cost.xlsx - Chris Transportation updates for AVMC
# (part of the broader "AVIC, AVLRaPTC, AVMC, BBSoEVP, BESP, BMRESP, BNVP,
# SoCDTtiNTY, TTS" transportation-variables update pass)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "About" sheet: insert a new note row right above the "Currency
#    Conversion" block (old row 74) explaining that the California model
#    re-uses the motorbike-freight category for long-haul class 8 trucks.
#    Inserting a real row shifts every row below it down by one, which
#    also re-points all the "About!$A$75" style formulas on the
#    "Cost Data" sheet to "About!$A$76", etc. automatically.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Rows.Item(73).Insert()
$wsAbout.Range("A73").Value = "The California model uses the motorbike-freight category to represent long-haul class 8 trucks. "
$wsAbout.Range("A73").ClearFormats()

# ---------------------------------------------------------------------
# 2. "Cost Data" sheet: the "Annual maintenance cost" row (row 100) used
#    to be an #N/A placeholder (=NA()) in columns B/C - point it at the
#    already-computed totals in row 96, and add an explanatory note in D.
# ---------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("Cost Data")
$wsCost.Range("B100").Formula = "=B96"
$wsCost.Range("C100").Formula = "=C96"
$wsCost.Range("D100").Value = "Note: The California model uses the motorbike-freight category to represent long-haul class 8 trucks. "
$wsCost.Range("D100").Font.ThemeColor = 1

# ---------------------------------------------------------------------
# 3. "AVMC-freight" sheet: row 7 ("typical annual maintenance cost")
#    used to be hard-coded zeros - have it follow row 3 like the other
#    metric rows on this sheet.
# ---------------------------------------------------------------------
$wsFreight = $wb.Worksheets.Item("AVMC-freight")
$wsFreight.Range("B7").Formula = "=B3"
$wsFreight.Range("C7:H7").Formula = "=C3"

# ---------------------------------------------------------------------
# 4. View/selection bookkeeping, matching where each sheet was left
#    scrolled/selected and which sheet ends up active.
# ---------------------------------------------------------------------
$wsAbout.Range("A73").Select()

$wsCost.Range("D101").Select()

$wsPassenger = $wb.Worksheets.Item("AVMC-passenger")
$wsPassenger.Range("B2").Select()

$wsFreight.Range("B7:H7").Select()

$wsPassenger.Activate()
